# Update countries & provincias Spain
#
# 1) Swap the Suecia/Mexico rows (Suecia now ranks above Mexico) and
#    refresh Suecia's daily figures.
# 2) Swap the Banglades/Dinamarca rows (Dinamarca now ranks above
#    Banglades) and refresh Dinamarca's daily figures.
# 3) Refresh the daily figures for several other countries whose stats
#    changed (Estados Unidos, Paises Bajos, Catar, Rumania, Uzbekistan,
#    Croacia) without moving their rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($rowNum, $values) {
    $ws.Cells.Item($rowNum, 1).Value = $values[0]
    $ws.Cells.Item($rowNum, 2).Value = $values[1]
    $ws.Cells.Item($rowNum, 3).Value = $values[2]
    $ws.Cells.Item($rowNum, 4).Value = $values[3]
    $ws.Cells.Item($rowNum, 5).Value = $values[4]
    $ws.Cells.Item($rowNum, 6).Value = $values[5]
    $ws.Cells.Item($rowNum, 7).Value = $values[6]
    $ws.Cells.Item($rowNum, 8).Value = $values[7]
}

# Row 4 - Estados Unidos (unchanged position, updated numbers)
Set-Row 4 @("Estados Unidos", 1160997, 223, 173725, 919824, 16475, 4, 67448)

# Row 18 - Paises Bajos (unchanged position, updated numbers)
Set-Row 18 @("Paises Bajos", 40571, 335, 0, 35265, 708, 69, 5056)

# Row 24 - now Suecia (was Mexico) with refreshed figures
Set-Row 24 @("Suecia", 22317, 235, 1005, 18633, 403, 10, 2679)

# Row 25 - now Mexico (was Suecia), carrying Mexico's previous figures
Set-Row 25 @("Mexico", 22088, 1349, 13447, 6580, 378, 89, 2061)

# Row 33 - Catar (unchanged position, updated numbers)
Set-Row 33 @("Catar", 15551, 679, 1664, 13875, 72, 0, 12)

# Row 37 - Rumania (unchanged position, updated numbers)
Set-Row 37 @("Rumania", 13163, 431, 4869, 7504, 255, 19, 790)

# Row 41 - now Dinamarca (was Banglades) with refreshed figures
Set-Row 41 @("Dinamarca", 9523, 116, 6987, 2052, 62, 9, 484)

# Row 42 - now Banglades (was Dinamarca), carrying Banglades' previous figures
Set-Row 42 @("Banglades", 9455, 665, 177, 9101, 1, 2, 177)

# Row 72 - Uzbekistan (unchanged position, updated numbers)
Set-Row 72 @("Uzbekistan", 2136, 18, 1308, 818, 8, 1, 10)

# Row 73 - Croacia (unchanged position, updated numbers)
Set-Row 73 @("Croacia", 2096, 8, 1489, 528, 15, 2, 79)
